$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037657214017382
$ws.Range("D2").Value = 1.03962185197878
$ws.Range("E2").Value = 1.045816212972785
$ws.Range("F2").Value = 1.055268759670722
$ws.Range("I2").Value = 1.039401280020186
$ws.Range("J2").Value = 1.042758894051824
$ws.Range("K2").Value = 1.042406517694913
$ws.Range("L2").Value = 1.048583392152967
$ws.Range("M2").Value = 1.058009683342602

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038489636446298
$ws.Range("D3").Value = 1.040226279985468
$ws.Range("E3").Value = 1.046567634715181
$ws.Range("F3").Value = 1.056127261425557
$ws.Range("I3").Value = 1.039584969961564
$ws.Range("J3").Value = 1.043236374408762
$ws.Range("K3").Value = 1.042821744567824
$ws.Range("L3").Value = 1.049146487206007
$ws.Range("M3").Value = 1.058681477803562

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039028952526641
$ws.Range("D4").Value = 1.040617925214927
$ws.Range("E4").Value = 1.047054865859398
$ws.Range("F4").Value = 1.056683970986044
$ws.Range("I4").Value = 1.039702999176619
$ws.Range("J4").Value = 1.043545347898728
$ws.Range("K4").Value = 1.043090267475393
$ws.Range("L4").Value = 1.049511191686241
$ws.Range("M4").Value = 1.059116729328961

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039255842842392
$ws.Range("D5").Value = 1.040782700347064
$ws.Range("E5").Value = 1.047259937991737
$ws.Range("F5").Value = 1.056918297341616
$ws.Range("I5").Value = 1.039752419123312
$ws.Range("J5").Value = 1.043675241870247
$ws.Range("K5").Value = 1.043203115753319
$ws.Range("L5").Value = 1.049664594404024
$ws.Range("M5").Value = 1.0592998403987

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039293948165491
$ws.Range("D6").Value = 1.040810374240537
$ws.Range("E6").Value = 1.0472943845334
$ws.Range("F6").Value = 1.056957658460994
$ws.Range("I6").Value = 1.039760705228033
$ws.Range("J6").Value = 1.043697051690712
$ws.Range("K6").Value = 1.043222061166899
$ws.Range("L6").Value = 1.049690356099561
$ws.Range("M6").Value = 1.05933059320536

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039031983612368
$ws.Range("D7").Value = 1.04062012644964
$ws.Range("E7").Value = 1.047057605101972
$ws.Range("F7").Value = 1.056687100946338
$ws.Range("I7").Value = 1.039703660313302
$ws.Range("J7").Value = 1.043547083544006
$ws.Range("K7").Value = 1.043091775513446
$ws.Range("L7").Value = 1.049513241144624
$ws.Range("M7").Value = 1.059119175554113

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037938392608946
$ws.Range("D8").Value = 1.039826008399547
$ws.Range("E8").Value = 1.046069948995897
$ws.Range("F8").Value = 1.055558644379338
$ws.Range("I8").Value = 1.039463530413806
$ws.Range("J8").Value = 1.042920257554019
$ws.Range("K8").Value = 1.042546877242584
$ws.Range("L8").Value = 1.048773620578193
$ws.Range("M8").Value = 1.05823660337743

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036016652681914
$ws.Range("D9").Value = 1.038430884260314
$ws.Range("E9").Value = 1.044337399349817
$ws.Range("F9").Value = 1.053579449988949
$ws.Range("I9").Value = 1.039034062500034
$ws.Range("J9").Value = 1.041815857817631
$ws.Range("K9").Value = 1.041585552519121
$ws.Range("L9").Value = 1.04747301719504
$ws.Range("M9").Value = 1.056685727941889

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0347391625163
$ws.Range("D10").Value = 1.037503740256265
$ws.Range("E10").Value = 1.04318774097839
$ws.Range("F10").Value = 1.052266352246686
$ws.Range("I10").Value = 1.038743538928407
$ws.Range("J10").Value = 1.041079768402846
$ws.Range("K10").Value = 1.040943974733442
$ws.Range("L10").Value = 1.046607853648552
$ws.Range("M10").Value = 1.055654822324959

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034186886995081
$ws.Range("D11").Value = 1.037102996137232
$ws.Range("E11").Value = 1.042691224849923
$ws.Range("F11").Value = 1.051699301326371
$ws.Range("I11").Value = 1.038616749192022
$ws.Range("J11").Value = 1.040761091981188
$ws.Range("K11").Value = 1.040666015512842
$ws.Range("L11").Value = 1.046233698886381
$ws.Range("M11").Value = 1.055209164296703

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033981882132093
$ws.Range("D12").Value = 1.036954251225275
$ws.Range("E12").Value = 1.042506992928689
$ws.Range("F12").Value = 1.051488905136042
$ws.Range("I12").Value = 1.038569505658709
$ws.Range("J12").Value = 1.040642730978104
$ws.Range("K12").Value = 1.040562747548167
$ws.Range("L12").Value = 1.046094792784868
$ws.Range("M12").Value = 1.055043738591739

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034025850240929
$ws.Range("D13").Value = 1.036986152538651
$ws.Range("E13").Value = 1.042546502388172
$ws.Range("F13").Value = 1.051534025312192
$ws.Range("I13").Value = 1.038579646250792
$ws.Range("J13").Value = 1.040668119358979
$ws.Range("K13").Value = 1.040584899831161
$ws.Range("L13").Value = 1.046124585349218
$ws.Range("M13").Value = 1.055079217892445

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034169938468171
$ws.Range("D14").Value = 1.037090698596249
$ws.Range("E14").Value = 1.042675992166666
$ws.Range("F14").Value = 1.051681905171796
$ws.Range("I14").Value = 1.038612847045454
$ws.Range("J14").Value = 1.040751308020078
$ws.Range("K14").Value = 1.040657479782593
$ws.Range("L14").Value = 1.046222215394603
$ws.Range("M14").Value = 1.055195487871212

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034258733901729
$ws.Range("D15").Value = 1.037155127410193
$ws.Range("E15").Value = 1.042755801159342
$ws.Range("F15").Value = 1.051773049606802
$ws.Range("I15").Value = 1.038633283537212
$ws.Range("J15").Value = 1.040802564612618
$ws.Range("K15").Value = 1.040702195870726
$ws.Range("L15").Value = 1.046282378025463
$ws.Range("M15").Value = 1.055267140472716

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034775834040809
$ws.Range("D16").Value = 1.037530351550986
$ws.Range("E16").Value = 1.043220720537265
$ws.Range("F16").Value = 1.052304017938913
$ws.Range("I16").Value = 1.038751932725752
$ws.Range("J16").Value = 1.041100919182329
$ws.Range("K16").Value = 1.040962418881775
$ws.Range("L16").Value = 1.046632695030962
$ws.Range("M16").Value = 1.055684414726656

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035100435758285
$ws.Range("D17").Value = 1.037765912549448
$ws.Range("E17").Value = 1.043512699863635
$ws.Range("F17").Value = 1.052637491163382
$ws.Range("I17").Value = 1.038826093388207
$ws.Range("J17").Value = 1.041288084821106
$ws.Range("K17").Value = 1.041125610157616
$ws.Range("L17").Value = 1.046852565665231
$ws.Range("M17").Value = 1.055946356739853

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035289855833303
$ws.Range("D18").Value = 1.037903380193526
$ws.Range("E18").Value = 1.043683131142098
$ws.Range("F18").Value = 1.052832147848979
$ws.Range("I18").Value = 1.03886925437128
$ws.Range("J18").Value = 1.041397260604269
$ws.Range("K18").Value = 1.041220782135417
$ws.Range("L18").Value = 1.046980857420805
$ws.Range("M18").Value = 1.056099213491026

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035354457614293
$ws.Range("D19").Value = 1.03795026473864
$ws.Range("E19").Value = 1.043741264892185
$ws.Range("F19").Value = 1.052898545687065
$ws.Range("I19").Value = 1.038883954900726
$ws.Range("J19").Value = 1.04143448757228
$ws.Range("K19").Value = 1.041253230782778
$ws.Range("L19").Value = 1.047024609169818
$ws.Range("M19").Value = 1.056151345582888

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035065600221808
$ws.Range("D20").Value = 1.037740631943937
$ws.Range("E20").Value = 1.04348136033525
$ws.Range("F20").Value = 1.052601697381786
$ws.Range("I20").Value = 1.038818146542044
$ws.Range("J20").Value = 1.041268003168603
$ws.Range("K20").Value = 1.041108102792808
$ws.Range("L20").Value = 1.046828970968327
$ws.Range("M20").Value = 1.055918245545951

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034127504340033
$ws.Range("D21").Value = 1.037059909356268
$ws.Range("E21").Value = 1.042637855220606
$ws.Range("F21").Value = 1.05163835182456
$ws.Range("I21").Value = 1.038603074323604
$ws.Range("J21").Value = 1.040726810756925
$ws.Range("K21").Value = 1.040636107380124
$ws.Range("L21").Value = 1.046193463789955
$ws.Range("M21").Value = 1.055161246167844

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033538467314367
$ws.Range("D22").Value = 1.036632545306829
$ws.Range("E22").Value = 1.042108646628636
$ws.Range("F22").Value = 1.051034000196324
$ws.Range("I22").Value = 1.038466992920038
$ws.Range("J22").Value = 1.040386597857716
$ws.Range("K22").Value = 1.040339221171757
$ws.Range("L22").Value = 1.045794310417573
$ws.Range("M22").Value = 1.054685936741194

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033850652372672
$ws.Range("D23").Value = 1.036859038441648
$ws.Range("E23").Value = 1.042389081737478
$ws.Range("F23").Value = 1.051354250544225
$ws.Range("I23").Value = 1.038539213234935
$ws.Range("J23").Value = 1.040566945378897
$ws.Range("K23").Value = 1.040496617492864
$ws.Range("L23").Value = 1.046005869310612
$ws.Range("M23").Value = 1.054937845441781

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035081340637651
$ws.Range("D24").Value = 1.037752054948731
$ws.Range("E24").Value = 1.043495520934157
$ws.Range("F24").Value = 1.052617870595802
$ws.Range("I24").Value = 1.038821737675473
$ws.Range("J24").Value = 1.041277077186965
$ws.Range("K24").Value = 1.041116013663053
$ws.Range("L24").Value = 1.046839632257591
$ws.Range("M24").Value = 1.05593094756712

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036512829141361
$ws.Range("D25").Value = 1.038791047064499
$ws.Range("E25").Value = 1.044784365868615
$ws.Range("F25").Value = 1.054090005983509
$ws.Range("I25").Value = 1.039145835748610
$ws.Range("J25").Value = 1.042101346105099
$ws.Range("K25").Value = 1.041834205911692
$ws.Range("L25").Value = 1.047808924977457
$ws.Range("M25").Value = 1.057086143002653
